# Fixed req ID #
# The "ID" column (A) for the Game Manager requirement rows (35-43) was
# mislabeled (it skipped GM_3). Relabel each row's ID down by one slot,
# restoring the missing "GM_3" at the top and dropping the duplicate-ish
# trailing "GM_12".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A35").Value = "GM_3"
$ws.Range("A36").Value = "GM_4"
$ws.Range("A37").Value = "GM_5"
$ws.Range("A38").Value = "GM_6"
$ws.Range("A39").Value = "GM_7"
$ws.Range("A40").Value = "GM_8"
$ws.Range("A41").Value = "GM_9"
$ws.Range("A42").Value = "GM_10"
$ws.Range("A43").Value = "GM_11"

# Restore the view state to match where the edit left the selection/scroll.
$ws.Activate()
$ws.Range("D43").Select()
